# Add new "I0" and "IF" columns (I:J) to the sheet, mirroring the existing
# header style from H1, and fill the data rows (2-64) with the matching
# values. In this dataset I0 and IF always carry the same number per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy formatting from H1 (bold, centered, bordered)
# onto I1:J1, then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-64): values per row, identical for I and J.
$values = @(9,7,9,9,8,8,8,7,9,7,8,7,8,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,11,7,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,5,6,6,5,4,3)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $v = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}

Write-Output "I0/IF columns populated"
